$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (C) column for all existing data rows (2-409)
# from 45205 to 45206.
for ($r = 2; $r -le 409; $r++) {
    $ws.Cells.Item($r, 3).Value2 = 45206
}

# Row 409 picks up an explicit row height (matches the other data rows).
$ws.Rows.Item(409).RowHeight = 15

# Append the new record as row 410.
$ws.Cells.Item(410, 1).Value = "A 47940-2023"

$ws.Cells.Item(410, 2).Value2 = 45204
$ws.Cells.Item(410, 2).NumberFormat = "YYYY-MM-DD"

$ws.Cells.Item(410, 3).Value2 = 45206
$ws.Cells.Item(410, 3).NumberFormat = "YYYY-MM-DD"

$ws.Cells.Item(410, 4).Value = "DALARNAS LÄN"
$ws.Cells.Item(410, 5).Value = "LUDVIKA"
$ws.Cells.Item(410, 6).Value = "Bergvik skog väst AB"
$ws.Cells.Item(410, 7).Value = 2.1
$ws.Cells.Item(410, 8).Value = 0
$ws.Cells.Item(410, 9).Value = 0
$ws.Cells.Item(410, 10).Value = 0
$ws.Cells.Item(410, 11).Value = 0
$ws.Cells.Item(410, 12).Value = 0
$ws.Cells.Item(410, 13).Value = 0
$ws.Cells.Item(410, 14).Value = 0
$ws.Cells.Item(410, 15).Value = 0
$ws.Cells.Item(410, 16).Value = 0
$ws.Cells.Item(410, 17).Value = 0

# Column R retains the wrap-text style even though it stays empty, like the
# row above it.
$ws.Cells.Item(410, 18).WrapText = $true
